$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3879292809067465
$ws.Range("A2").Value = -0.0099999996375963462
$ws.Range("A3").Value = -0.0089999996319587439
$ws.Range("A4").Value = -0.011999999898247182
$ws.Range("A5").Value = -0.0059999996389192845
$ws.Range("A6").Value = -0.0059999996286599355
$ws.Range("A7").Value = -0.019999999555022185
$ws.Range("A8").Value = -0.019999999554380921
$ws.Range("A9").Value = -0.0059999996281314694
$ws.Range("A10").Value = 0.052184568322303448
$ws.Range("A11").Value = -0.0044999996371295481
$ws.Range("A12").Value = -0.0059999996283166546
$ws.Range("A13").Value = -0.0059999996246933307
$ws.Range("A14").Value = -0.011999999592353205
$ws.Range("A15").Value = -0.0059999996232198427
$ws.Range("A16").Value = -0.0059999996229911368
$ws.Range("A17").Value = -0.0059999996231008268
$ws.Range("A18").Value = 0.0066422995748594005
$ws.Range("A19").Value = -0.0089999996415905947
$ws.Range("A20").Value = -0.085691840534604324
$ws.Range("A21").Value = -0.0089999996274503502
$ws.Range("A22").Value = -0.0089999996270750948
$ws.Range("A23").Value = -0.0089999996286298511
$ws.Range("A24").Value = -0.041999999450220038
$ws.Range("A25").Value = -0.041999999447278391
$ws.Range("A26").Value = -0.0059999996272601663
$ws.Range("A27").Value = -0.0059999996252790844
$ws.Range("A28").Value = -0.0059999996170194692
$ws.Range("A29").Value = -0.011999999580252663
$ws.Range("A30").Value = -0.019999999535726065
$ws.Range("A31").Value = -0.014999999557931076
$ws.Range("A32").Value = -0.02099999952585474
$ws.Range("A33").Value = -0.0059999996041142367
